# Update gh-pages output data (F column "浏览/收藏" counters) for the
# "展览" and "全部类型" sheets, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for sheet "展览" (sheet1.xml)
$updatesExhibition = @{
    2  = 629
    5  = 13136
    6  = 74
    8  = 517
    9  = 480
    10 = 1181
    11 = 988
    12 = 13769
    13 = 14369
    15 = 173
    21 = 34
    22 = 1093
    25 = 5420
    26 = 938
    27 = 21
    28 = 316
    30 = 46
}

# Map of row -> new F value for sheet "全部类型" (sheet4.xml)
$updatesAll = @{
    2  = 629
    5  = 13137
    6  = 74
    8  = 517
    9  = 480
    10 = 1181
    11 = 988
    12 = 13769
    13 = 14369
    15 = 173
    21 = 34
    22 = 1093
    25 = 5420
    26 = 938
    27 = 21
    28 = 316
    30 = 46
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $updatesExhibition[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $wsAll.Range("F$row").Value = $updatesAll[$row]
}
